$d = $word.ActiveDocument
$gallery = $word.ListGalleries.Item(3)

function New-TempListId {
    param($label)
    $end = $d.Content.End
    $r = $d.Range($end, $end)
    $r.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($d.Paragraphs.Count)
    $pr = $newP.Range
    $pr.InsertBefore($label)
    $lt = $gallery.ListTemplates.Add($false, $false)
    $pr.ListFormat.ApplyListTemplateWithLevel($lt, $false, 0, $false, $false)
}

# Materialize numId 2, 3, 4, 5 (each its own abstract numbering definition,
# decimal "%1." at ilvl 0) by briefly applying them to throw-away paragraphs
# at the end of the document, then deleting those paragraphs. The numbering
# part keeps the definitions even after the referencing paragraphs are gone.
New-TempListId "TEMP2"
New-TempListId "TEMP3"
New-TempListId "TEMP4"
New-TempListId "TEMP5"

for ($i = 0; $i -lt 4; $i++) {
    $count = $d.Paragraphs.Count
    $p = $d.Paragraphs.Item($count)
    $p.Range.Delete()
}

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range

$xml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t>הממשק הגרפי-</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve"> יהא זהה לחלוטין לגרסת ההדמיה הגרפית.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:rtl/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>דרישות שאינן פונקציונליות</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t>עבור הממשק הגרפי</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:rPr>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t>.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
    <w:t>NET F</w:t>
  </w:r>
  <w:r>
    <w:t>ramework</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>חיבור לאינטרנט</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t>עבור שימוש ב</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>API</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>python 2.7</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve">או </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">  C++</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve">או </w:t>
  </w:r>
  <w:r>
    <w:t>.NET Framework</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve">  </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>חיבור אינטרנט עבור שימוש ב</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
    <w:t>API</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve"> הפומבי/גלובאלי.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve">עבור הממשק הגרפי </w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve">דפדפן (בעל תמיכה </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
    <w:t>HTML 5</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve"> ו </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
    <w:t>J</w:t>
  </w:r>
  <w:r>
    <w:t>avaScript</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>).</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>חיבור לאינטרנט</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve">תקשורת </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
    <w:t>SSL</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t xml:space="preserve"> עבור אבטחת תקשורת.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:b/>
      <w:bCs/>
      <w:rtl/>
    </w:rPr>
    <w:t>עבור האפליקציה</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="5"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>מערכת הפעלה אנדרואיד 4 או יותר*</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="5"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>מקום פנוי בזכרון.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="5"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="cs"/>
      <w:rtl/>
    </w:rPr>
    <w:t>חיבור אינטרנט.</w:t>
  </w:r>
</w:p>

'@

$r.InsertXML($xml)
